$wb = $excel.ActiveWorkbook
Get-Member -InputObject $wb
